$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range("D2").Style
$ws.Range("D2").Formula = "'94.307.47"
$ws.Range("D2").Style = $style
$style = $ws.Range("E2").Style
$ws.Range("E2").Formula = "'  +2.55%  "
$ws.Range("E2").Style = $style

# Row 3
$style = $ws.Range("D3").Style
$ws.Range("D3").Formula = "'3.099.74"
$ws.Range("D3").Style = $style
$style = $ws.Range("E3").Style
$ws.Range("E3").Formula = "'  +0.01%  "
$ws.Range("E3").Style = $style

# Row 4
$style = $ws.Range("E4").Style
$ws.Range("E4").Formula = "'  +0.02%  "
$ws.Range("E4").Style = $style

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").Formula = "'236.92"
$ws.Range("D5").Style = $style
$style = $ws.Range("E5").Style
$ws.Range("E5").Formula = "'  -1.05%  "
$ws.Range("E5").Style = $style

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").Formula = "'613.02"
$ws.Range("D6").Style = $style
$style = $ws.Range("E6").Style
$ws.Range("E6").Formula = "'  -0.12%  "
$ws.Range("E6").Style = $style

# Row 7
$style = $ws.Range("E7").Style
$ws.Range("E7").Formula = "'  +3.60%  "
$ws.Range("E7").Style = $style

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").Formula = "'0.391"
$ws.Range("D8").Style = $style
$style = $ws.Range("E8").Style
$ws.Range("E8").Formula = "'  +0.21%  "
$ws.Range("E8").Style = $style

# Row 9
$style = $ws.Range("E9").Style
$ws.Range("E9").Formula = "'  -0.02%  "
$ws.Range("E9").Style = $style

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").Formula = "'0.821"
$ws.Range("D10").Style = $style
$style = $ws.Range("E10").Style
$ws.Range("E10").Formula = "'  +12.84%  "
$ws.Range("E10").Style = $style

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").Formula = "'3.102.61"
$ws.Range("D11").Style = $style
$style = $ws.Range("E11").Style
$ws.Range("E11").Formula = "'  +0.11%  "
$ws.Range("E11").Style = $style

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").Formula = "'0.198"
$ws.Range("D12").Style = $style
$style = $ws.Range("E12").Style
$ws.Range("E12").Formula = "'  -1.78%  "
$ws.Range("E12").Style = $style

# Row 13
$style = $ws.Range("B13").Style
$ws.Range("B13").Formula = "'WrappedBTC"
$ws.Range("B13").Style = $style
$style = $ws.Range("C13").Style
$ws.Range("C13").Formula = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C13").Style = $style
$style = $ws.Range("D13").Style
$ws.Range("D13").Formula = "'94.095.87"
$ws.Range("D13").Style = $style
$style = $ws.Range("E13").Style
$ws.Range("E13").Formula = "'  +2.41%  "
$ws.Range("E13").Style = $style

# Row 14
$style = $ws.Range("B14").Style
$ws.Range("B14").Formula = "'ShibaInu"
$ws.Range("B14").Style = $style
$style = $ws.Range("C14").Style
$ws.Range("C14").Formula = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C14").Style = $style
$style = $ws.Range("D14").Style
$ws.Range("D14").Formula = "'0.0000244"
$ws.Range("D14").Style = $style
$style = $ws.Range("E14").Style
$ws.Range("E14").Formula = "'  -2.34%  "
$ws.Range("E14").Style = $style

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").Formula = "'34.43"
$ws.Range("D15").Style = $style
$style = $ws.Range("E15").Style
$ws.Range("E15").Formula = "'  +0.35%  "
$ws.Range("E15").Style = $style

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").Formula = "'5.37"
$ws.Range("D16").Style = $style
$style = $ws.Range("E16").Style
$ws.Range("E16").Formula = "'  -2.42%  "
$ws.Range("E16").Style = $style

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").Formula = "'3.689.57"
$ws.Range("D17").Style = $style
$style = $ws.Range("E17").Style
$ws.Range("E17").Formula = "'  +0.12%  "
$ws.Range("E17").Style = $style

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").Formula = "'3.122.60"
$ws.Range("D18").Style = $style
$style = $ws.Range("E18").Style
$ws.Range("E18").Formula = "'  +1.10%  "
$ws.Range("E18").Style = $style

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").Formula = "'3.68"
$ws.Range("D19").Style = $style
$style = $ws.Range("E19").Style
$ws.Range("E19").Formula = "'  +1.31%  "
$ws.Range("E19").Style = $style

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").Formula = "'14.82"
$ws.Range("D20").Style = $style
$style = $ws.Range("E20").Style
$ws.Range("E20").Formula = "'  +0.74%  "
$ws.Range("E20").Style = $style

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").Formula = "'5.90"
$ws.Range("D21").Style = $style
$style = $ws.Range("E21").Style
$ws.Range("E21").Formula = "'  +1.81%  "
$ws.Range("E21").Style = $style

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").Formula = "'448.79"
$ws.Range("D22").Style = $style
$style = $ws.Range("E22").Style
$ws.Range("E22").Formula = "'  +0.66%  "
$ws.Range("E22").Style = $style

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").Formula = "'0.0000198"
$ws.Range("D23").Style = $style
$style = $ws.Range("E23").Style
$ws.Range("E23").Formula = "'  -1.40%  "
$ws.Range("E23").Style = $style

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").Formula = "'8.94"
$ws.Range("D24").Style = $style
$style = $ws.Range("E24").Style
$ws.Range("E24").Formula = "'  -3.82%  "
$ws.Range("E24").Style = $style

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").Formula = "'8.30"
$ws.Range("D25").Style = $style
$style = $ws.Range("E25").Style
$ws.Range("E25").Formula = "'  +5.52%  "
$ws.Range("E25").Style = $style

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").Formula = "'5.62"
$ws.Range("D26").Style = $style
$style = $ws.Range("E26").Style
$ws.Range("E26").Formula = "'  +0.32%  "
$ws.Range("E26").Style = $style

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").Formula = "'86.19"
$ws.Range("D27").Style = $style
$style = $ws.Range("E27").Style
$ws.Range("E27").Formula = "'  -0.79%  "
$ws.Range("E27").Style = $style

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").Formula = "'12.01"
$ws.Range("D28").Style = $style
$style = $ws.Range("E28").Style
$ws.Range("E28").Formula = "'  +3.19%  "
$ws.Range("E28").Style = $style

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").Formula = "'3.284.15"
$ws.Range("D29").Style = $style
$style = $ws.Range("E29").Style
$ws.Range("E29").Formula = "'  +0.26%  "
$ws.Range("E29").Style = $style

# Row 30
$style = $ws.Range("E30").Style
$ws.Range("E30").Formula = "'  +0.15%  "
$ws.Range("E30").Style = $style

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").Formula = "'0.254"
$ws.Range("D31").Style = $style
$style = $ws.Range("E31").Style
$ws.Range("E31").Formula = "'  +9.78%  "
$ws.Range("E31").Style = $style

# Row 32
$style = $ws.Range("E32").Style
$ws.Range("E32").Formula = "'  +8.27%  "
$ws.Range("E32").Style = $style

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").Formula = "'0.126"
$ws.Range("D33").Style = $style
$style = $ws.Range("E33").Style
$ws.Range("E33").Formula = "'  -7.59%  "
$ws.Range("E33").Style = $style

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").Formula = "'9.27"
$ws.Range("D34").Style = $style
$style = $ws.Range("E34").Style
$ws.Range("E34").Formula = "'  +0.81%  "
$ws.Range("E34").Style = $style

# Row 35
$style = $ws.Range("E35").Style
$ws.Range("E35").Formula = "'  +0.31%  "
$ws.Range("E35").Style = $style

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").Formula = "'7.86"
$ws.Range("D36").Style = $style
$style = $ws.Range("E36").Style
$ws.Range("E36").Formula = "'  -0.59%  "
$ws.Range("E36").Style = $style

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").Formula = "'0.161"
$ws.Range("D37").Style = $style
$style = $ws.Range("E37").Style
$ws.Range("E37").Formula = "'  -2.70%  "
$ws.Range("E37").Style = $style

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").Formula = "'25.98"
$ws.Range("D38").Style = $style
$style = $ws.Range("E38").Style
$ws.Range("E38").Formula = "'  -0.39%  "
$ws.Range("E38").Style = $style

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").Formula = "'1.91"
$ws.Range("D39").Style = $style
$style = $ws.Range("E39").Style
$ws.Range("E39").Formula = "'  -1.38%  "
$ws.Range("E39").Style = $style

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").Formula = "'0.452"
$ws.Range("D40").Style = $style
$style = $ws.Range("E40").Style
$ws.Range("E40").Formula = "'  +5.12%  "
$ws.Range("E40").Style = $style

# Row 41
$style = $ws.Range("B41").Style
$ws.Range("B41").Formula = "'WhiteBITCoin"
$ws.Range("B41").Style = $style
$style = $ws.Range("C41").Style
$ws.Range("C41").Formula = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C41").Style = $style
$style = $ws.Range("D41").Style
$ws.Range("D41").Formula = "'23.95"
$ws.Range("D41").Style = $style
$style = $ws.Range("E41").Style
$ws.Range("E41").Formula = "'  +7.98%  "
$ws.Range("E41").Style = $style

# Row 42
$style = $ws.Range("B42").Style
$ws.Range("B42").Formula = "'Bittensor"
$ws.Range("B42").Style = $style
$style = $ws.Range("C42").Style
$ws.Range("C42").Formula = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C42").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").Formula = "'473.83"
$ws.Range("D42").Style = $style
$style = $ws.Range("E42").Style
$ws.Range("E42").Formula = "'  -1.29%  "
$ws.Range("E42").Style = $style

# Row 43
$style = $ws.Range("E43").Style
$ws.Range("E43").Formula = "'  -1.16%  "
$ws.Range("E43").Style = $style

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").Formula = "'3.69"
$ws.Range("D44").Style = $style
$style = $ws.Range("E44").Style
$ws.Range("E44").Formula = "'  -12.69%  "
$ws.Range("E44").Style = $style

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").Formula = "'3.27"
$ws.Range("D45").Style = $style
$style = $ws.Range("E45").Style
$ws.Range("E45").Formula = "'  -4.08%  "
$ws.Range("E45").Style = $style

# Row 46
$style = $ws.Range("E46").Style
$ws.Range("E46").Formula = "'  +0.00%  "
$ws.Range("E46").Style = $style

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").Formula = "'160.53"
$ws.Range("D47").Style = $style
$style = $ws.Range("E47").Style
$ws.Range("E47").Formula = "'  +1.06%  "
$ws.Range("E47").Style = $style

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").Formula = "'0.687"
$ws.Range("D48").Style = $style
$style = $ws.Range("E48").Style
$ws.Range("E48").Formula = "'  -0.85%  "
$ws.Range("E48").Style = $style

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").Formula = "'1.85"
$ws.Range("D49").Style = $style
$style = $ws.Range("E49").Style
$ws.Range("E49").Formula = "'  -2.23%  "
$ws.Range("E49").Style = $style

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").Formula = "'4.41"
$ws.Range("D50").Style = $style
$style = $ws.Range("E50").Style
$ws.Range("E50").Formula = "'  +0.50%  "
$ws.Range("E50").Style = $style

# Row 51
$style = $ws.Range("B51").Style
$ws.Range("B51").Formula = "'OKB"
$ws.Range("B51").Style = $style
$style = $ws.Range("C51").Style
$ws.Range("C51").Formula = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C51").Style = $style
$style = $ws.Range("D51").Style
$ws.Range("D51").Formula = "'43.82"
$ws.Range("D51").Style = $style
$style = $ws.Range("E51").Style
$ws.Range("E51").Formula = "'  -0.45%  "
$ws.Range("E51").Style = $style
